$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.160498738288879
$ws.Range("B1").Value = 2.410387516021729
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.383260011672974
$ws.Range("E1").Value = 1.230436444282532
